$wb = $excel.ActiveWorkbook

# --- Segment Series sheet: add a "Time Bonus (secs)" column (I) ---
$ws = $wb.Worksheets.Item("Segment Series")
$ws.Activate()

# Copy the existing header formatting from H1 onto the new I1 header cell
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "Time Bonus (secs)"
$ws.Range("I3").Value = 60
$ws.Range("I4").Value = -30

# Size the new column to fit its header text
$ws.Columns.Item(9).ColumnWidth = 15.5

# Leave the selection on the newly-added data
$ws.Range("I10:J10").Select() | Out-Null
